# Create Requirements - Item . Link to supplier flow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), written left-to-right so the shared-string
#     table gets the new entries appended in this exact order ---
$ws.Range("F1").Value = "FormTask"
$ws.Range("G1").Value = "DocTask"
$ws.Range("H1").Value = "ACKTask"
$ws.Range("I1").Value = "ItemName"
$ws.Range("J1").Value = "FormTask1"
$ws.Range("K1").Value = "DocTask1"
$ws.Range("L1").Value = "ACKTask1"

# --- New / updated data cells (row 2), in authoring order ---
$ws.Range("H2").Value = "AUTO_TEST_RESTORE_DOCUMENT_ON_03/01/2019-12:40:27"
$ws.Range("A2").Value = "Supplier_On_11/01/19-14:15"
$ws.Range("E2").Value = "WorkGroup_11/01/19-14:15"
$ws.Range("D2").Value = "Req_Supplier_On_11/01/19-14:16"
$ws.Range("F2").Value = "TestFormTask_11/01/19-14:17"
$ws.Range("I2").Value = "Req_Supplier_On_11/01/19-14:18"
$ws.Range("L2").Value = "TestACKTask_11/01/19-14:19"
$ws.Range("J2").Value = "TestFormTask_11/01/19-14:20"
$ws.Range("K2").Value = "TestDocTask_11/01/19-14:21"
$ws.Range("G2").Value = "jj"

# --- Column widths for the two new columns ---
$ws.Columns.Item(9).ColumnWidth = 16
$ws.Columns.Item(10).ColumnWidth = 24.833333333333336

# --- Move the active selection ---
$ws.Range("J5").Select()
